# Fruta / hortaliza, semanal
# Insert a new price-record row at row 196 (pushing the old rows 196-198
# down to 197-199), and populate it with the new weekly record for
# "Florida King" durazno from Comercializadora del Agro de Limarí.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above the current row 196; this shifts the existing
# rows 196, 197, 198 down to 197, 198, 199 and carries formatting
# (including the date style on column D) down with them.
$ws.Rows(196).Insert()

# Populate the newly inserted row 196 with the new record.
$ws.Cells.Item(196, 1).Value  = 2
$ws.Cells.Item(196, 2).Value  = "Comercializadora del Agro de Limarí"
$ws.Cells.Item(196, 3).Value  = "Coquimbo"
$ws.Cells.Item(196, 4).Value  = 45239
$ws.Cells.Item(196, 5).Value  = 4
$ws.Cells.Item(196, 6).Value  = "Fruta"
$ws.Cells.Item(196, 7).Value  = 100103
$ws.Cells.Item(196, 8).Value  = "Frutos de hueso (carozo)"
$ws.Cells.Item(196, 9).Value  = 100103004
$ws.Cells.Item(196, 10).Value = "Durazno"
$ws.Cells.Item(196, 11).Value = "Florida King"
$ws.Cells.Item(196, 12).Value = "Primera"
$ws.Cells.Item(196, 13).Value = 300
$ws.Cells.Item(196, 14).Value = 14000
$ws.Cells.Item(196, 15).Value = 15000
$ws.Cells.Item(196, 16).Value = 14500
$ws.Cells.Item(196, 17).Value = "$/bandeja 10 kilos granel"
$ws.Cells.Item(196, 18).Value = "Provincia de Limarí"
$ws.Cells.Item(196, 19).Value = 1450
$ws.Cells.Item(196, 20).Value = 10
